$d = $word.ActiveDocument

# UCM+SP row: 93 (30.4) -> 94 (30.6)
$d.Content.Find.Execute("93 (30.4)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "94 (30.6)", 2)

# UCM row: 94 (30.7) -> 94 (30.6)
$d.Content.Find.Execute("94 (30.7)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "94 (30.6)", 2)

# L2C row: 91 (29.7) -> 91 (29.6)
$d.Content.Find.Execute("91 (29.7)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "91 (29.6)", 2)

# Total row: 306 (100.0) -> 307 (100.0)
$d.Content.Find.Execute("306 (100.0)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "307 (100.0)", 2)
